# Update VIF data table: revised variable labels/values and 2 fewer rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (row 1 header "Variable"/"VIF" unchanged)
$data = @(
    @("Sex", 1.055649933362003),
    @("Age", 1.069351908737812),
    @("Motor function", 1.332338024332296),
    @("Dermatological symptoms", 1.08127001191637),
    @("Arthralgia", 1.134161850216428),
    @("Urinary function impairment", 1.63706327914107),
    @("Lower limb pain", 1.171526008238892),
    @("Paresthesia", 1.165442856382306),
    @("Lower back spine pain", 1.314918936897601),
    @("Limbs paresis or weakness", 1.402443153044091),
    @("CD39-Diplotypes", 1.033543457305548)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-obsolete rows 13 and 14 (table shrank from 14 to 12 rows)
$ws.Range("A13:B14").Clear()
